# Auto-generated edit script: applies cell value updates across 8 sheets
# matching the authoritative diff of Sheets/Midgardsormr_Profits.xlsx
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 742.1111
$ws.Range("I9").Value = 696.3333
$ws.Range("J9").Value = 765
$ws.Range("K9").Value = 696.3333
$ws.Range("L9").Value = 765
$ws.Range("M9").Value = -527.3333
$ws.Range("N9").Value = -1103
$ws.Range("H43").Value = 5590.6665
$ws.Range("I43").Value = 4939.6
$ws.Range("J43").Value = 6055.7144
$ws.Range("K43").Value = 4939.6
$ws.Range("L43").Value = 6055.7144
$ws.Range("M43").Value = -4870.6
$ws.Range("N43").Value = -6193.7144
$ws.Range("H53").Value = 724.7646999999999
$ws.Range("I53").Value = 613.1429000000001
$ws.Range("J53").Value = 802.9
$ws.Range("K53").Value = 613.1429000000001
$ws.Range("L53").Value = 802.9
$ws.Range("M53").Value = 23.85709999999995
$ws.Range("N53").Value = -2076.9
$ws.Range("H107").Value = 1403.8857
$ws.Range("I107").Value = 1113.32
$ws.Range("J107").Value = 2130.3
$ws.Range("K107").Value = 1113.32
$ws.Range("L107").Value = 2130.3
$ws.Range("M107").Value = 806.6800000000001
$ws.Range("N107").Value = -5970.3
$ws.Range("H116").Value = 22537.408
$ws.Range("I116").Value = 19710
$ws.Range("J116").Value = 25930.3
$ws.Range("K116").Value = 19710
$ws.Range("L116").Value = 25930.3
$ws.Range("M116").Value = -16268
$ws.Range("N116").Value = -32814.3
$ws.Range("H137").Value = 41032.668
$ws.Range("I137").Value = 67999.2
$ws.Range("J137").Value = 7324.5
$ws.Range("K137").Value = 203997.6
$ws.Range("L137").Value = 21973.5
$ws.Range("M137").Value = -201447.6
$ws.Range("N137").Value = -27073.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 266.66666
$ws.Range("I17").Value = 266.66666
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 266.66666
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -93.66665999999998
$ws.Range("N17").ClearContents()
$ws.Range("H32").Value = 20347.9
$ws.Range("I32").Value = 14768.584
$ws.Range("K32").Value = 14768.584
$ws.Range("M32").Value = -14481.584
$ws.Range("H43").Value = 81079.75
$ws.Range("J43").Value = 99459
$ws.Range("L43").Value = 99459
$ws.Range("N43").Value = -100085
$ws.Range("H45").Value = 4609.4
$ws.Range("I45").Value = 2398.6
$ws.Range("K45").Value = 2398.6
$ws.Range("M45").Value = -2021.6
$ws.Range("H74").Value = 136449.16
$ws.Range("I74").Value = 177290.14
$ws.Range("J74").Value = 10213.363
$ws.Range("K74").Value = 177290.14
$ws.Range("L74").Value = 10213.363
$ws.Range("M74").Value = -176416.14
$ws.Range("N74").Value = -11961.363
$ws.Range("H77").Value = 136449.16
$ws.Range("I77").Value = 177290.14
$ws.Range("J77").Value = 10213.363
$ws.Range("K77").Value = 886450.7000000001
$ws.Range("L77").Value = 51066.815
$ws.Range("M77").Value = -882082.7000000001
$ws.Range("N77").Value = -59802.815
$ws.Range("H92").Value = 1027163.5
$ws.Range("I92").Value = 90000
$ws.Range("K92").Value = 90000
$ws.Range("M92").Value = -87504
$ws.Range("H132").Value = 1005.19446
$ws.Range("I132").Value = 961.58826
$ws.Range("J132").Value = 1746.5
$ws.Range("K132").Value = 2884.76478
$ws.Range("L132").Value = 5239.5
$ws.Range("M132").Value = -354.76478
$ws.Range("N132").Value = -10299.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1367.6857
$ws.Range("I86").Value = 1319.8276
$ws.Range("J86").Value = 1599
$ws.Range("K86").Value = 1319.8276
$ws.Range("L86").Value = 1599
$ws.Range("M86").Value = -196.8276000000001
$ws.Range("N86").Value = -3845
$ws.Range("H89").Value = 1367.6857
$ws.Range("I89").Value = 1319.8276
$ws.Range("J89").Value = 1599
$ws.Range("K89").Value = 6599.138000000001
$ws.Range("L89").Value = 7995
$ws.Range("M89").Value = -983.1380000000008
$ws.Range("N89").Value = -19227
$ws.Range("H94").Value = 4030.0625
$ws.Range("I94").Value = 5028.727
$ws.Range("J94").Value = 1833
$ws.Range("K94").Value = 5028.727
$ws.Range("L94").Value = 1833
$ws.Range("M94").Value = -4577.727
$ws.Range("N94").Value = -2735
$ws.Range("H105").Value = 2936.2188
$ws.Range("I105").Value = 2710.5
$ws.Range("J105").Value = 3914.3333
$ws.Range("K105").Value = 2710.5
$ws.Range("L105").Value = 3914.3333
$ws.Range("M105").Value = -963.5
$ws.Range("N105").Value = -7408.3333
$ws.Range("H134").Value = 3072.7273
$ws.Range("I134").Value = 1487.8334
$ws.Range("K134").Value = 4463.5002
$ws.Range("M134").Value = -1928.5002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10003898
$ws.Range("I31").Value = 20002398
$ws.Range("J31").Value = 5398.2
$ws.Range("K31").Value = 20002398
$ws.Range("L31").Value = 5398.2
$ws.Range("M31").Value = -20002103
$ws.Range("N31").Value = -5988.2
$ws.Range("H34").Value = 10003898
$ws.Range("I34").Value = 20002398
$ws.Range("J34").Value = 5398.2
$ws.Range("K34").Value = 20002398
$ws.Range("L34").Value = 5398.2
$ws.Range("M34").Value = -20002196
$ws.Range("N34").Value = -5802.2
$ws.Range("H135").Value = 97799.914
$ws.Range("J135").Value = 97799.914
$ws.Range("L135").Value = 97799.914
$ws.Range("N135").Value = -107939.914

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 199.28572
$ws.Range("I14").Value = 199.28572
$ws.Range("K14").Value = 597.85716
$ws.Range("M14").Value = -424.85716
$ws.Range("H107").Value = 3720.6
$ws.Range("I107").Value = 15201.5
$ws.Range("J107").Value = 850.375
$ws.Range("K107").Value = 45604.5
$ws.Range("L107").Value = 2551.125
$ws.Range("M107").Value = -43684.5
$ws.Range("N107").Value = -6391.125
$ws.Range("H129").Value = 2190.111
$ws.Range("I129").Value = 1572
$ws.Range("J129").Value = 2962.75
$ws.Range("K129").Value = 4716
$ws.Range("L129").Value = 8888.25
$ws.Range("M129").Value = 284
$ws.Range("N129").Value = -18888.25
$ws.Range("H132").Value = 2017.7142
$ws.Range("J132").Value = 1799
$ws.Range("L132").Value = 16191
$ws.Range("N132").Value = -21251
$ws.Range("H137").Value = 3846.3333
$ws.Range("I137").Value = 2910
$ws.Range("J137").Value = 4515.143
$ws.Range("K137").Value = 8730
$ws.Range("L137").Value = 13545.429
$ws.Range("M137").Value = -3630
$ws.Range("N137").Value = -23745.429
$ws.Range("H139").Value = 5134.706
$ws.Range("I139").Value = 5019.3335
$ws.Range("J139").Value = 6000
$ws.Range("K139").Value = 15058.0005
$ws.Range("L139").Value = 18000
$ws.Range("M139").Value = -9918.000499999998
$ws.Range("N139").Value = -28280

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1563.5454
$ws.Range("I97").Value = 1299.3529
$ws.Range("J97").Value = 2461.8
$ws.Range("K97").Value = 1299.3529
$ws.Range("L97").Value = 2461.8
$ws.Range("M97").Value = -803.3529000000001
$ws.Range("N97").Value = -3453.8
$ws.Range("H126").Value = 2628.1052
$ws.Range("I126").Value = 1785.0834
$ws.Range("K126").Value = 5355.2502
$ws.Range("M126").Value = -2885.2502
$ws.Range("H132").Value = 3736.2964
$ws.Range("I132").Value = 4030.3635
$ws.Range("K132").Value = 12091.0905
$ws.Range("M132").Value = -9561.0905
$ws.Range("H140").Value = 77852.28999999999
$ws.Range("J140").Value = 77852.28999999999
$ws.Range("L140").Value = 77852.28999999999
$ws.Range("N140").Value = -88212.28999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 7721.154
$ws.Range("I100").Value = 4330.5
$ws.Range("K100").Value = 4330.5
$ws.Range("M100").Value = -3789.5
$ws.Range("H122").Value = 6964
$ws.Range("I122").Value = 5474.75
$ws.Range("K122").Value = 16424.25
$ws.Range("M122").Value = -13974.25
$ws.Range("H132").Value = 3380.9092
$ws.Range("I132").Value = 3077.2222
$ws.Range("K132").Value = 9231.6666
$ws.Range("M132").Value = -6701.6666

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5856.3335
$ws.Range("I62").Value = 5697
$ws.Range("J62").Value = 5936
$ws.Range("K62").Value = 5697
$ws.Range("L62").Value = 5936
$ws.Range("M62").Value = -5073
$ws.Range("N62").Value = -7184
$ws.Range("H65").Value = 5856.3335
$ws.Range("I65").Value = 5697
$ws.Range("J65").Value = 5936
$ws.Range("K65").Value = 28485
$ws.Range("L65").Value = 29680
$ws.Range("M65").Value = -25365
$ws.Range("N65").Value = -35920
$ws.Range("H107").Value = 1124.4
$ws.Range("I107").Value = 867.1111
$ws.Range("J107").Value = 1334.909
$ws.Range("K107").Value = 2601.3333
$ws.Range("L107").Value = 4004.727
$ws.Range("M107").Value = -681.3332999999998
$ws.Range("N107").Value = -7844.727000000001
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 142107
$ws.Range("J131").Value = 142107
$ws.Range("L131").Value = 142107
$ws.Range("N131").Value = -152187
$ws.Range("H132").Value = 3116.4
$ws.Range("I132").Value = 2954.0938
$ws.Range("K132").Value = 8862.2814
$ws.Range("M132").Value = -6332.2814
$ws.Range("H137").Value = 110649.29
$ws.Range("J137").Value = 110649.29
$ws.Range("L137").Value = 110649.29
$ws.Range("N137").Value = -120849.29
